# Updates cryptos list values (Price and Volume(1h) columns) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.287.48'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '1.933.25'
$ws.Range('E3').Value = '  +1.56%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'325.38"
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').Value = "'0.9999"
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').Value = "'0.4625"
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').Value = "'0.3876"
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = "'45.87"
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('D10').Value = "'0.07818"
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('D11').Value = "'0.9757"
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('D12').Value = "'22.68"
$ws.Range('E12').Value = '  +3.32%  '
$ws.Range('D13').Value = '1.932.13'
$ws.Range('E13').Value = '  +3.62%  '
$ws.Range('D14').Value = "'5.787"
$ws.Range('E14').Value = '  +0.73%  '
$ws.Range('D15').Value = "'7.084"
$ws.Range('E15').Value = '  +0.75%  '
$ws.Range('D16').Value = "'0.07061"
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').Value = "'86.75"
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('D18').Value = "'1.002"
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('D19').Value = "'0.000009746"
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').Value = "'17.03"
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').Value = '29.274.33'
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('E23').Value = '  +2.83%  '
$ws.Range('D24').Value = "'11.06"
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('D25').Value = '2.160.93'
$ws.Range('E25').Value = '  +3.00%  '
$ws.Range('D26').Value = "'2.093"
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').Value = "'157.58"
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('D28').Value = "'19.36"
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('D29').Value = "'5.766"
$ws.Range('E29').Value = '  -2.16%  '
$ws.Range('D30').Value = "'118.88"
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('D31').Value = "'1.838"
$ws.Range('E31').Value = '  -1.95%  '
$ws.Range('D32').Value = "'0.09345"
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').Value = "'0.8644"
$ws.Range('E33').Value = '  -3.32%  '
$ws.Range('D34').Value = "'5.167"
$ws.Range('E34').Value = '  -1.06%  '
$ws.Range('D35').Value = "'1.301"
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('E36').Value = '  -1.76%  '
$ws.Range('D37').Value = "'0.05782"
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').Value = "'1.155"
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').Value = "'0.02080"
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').Value = "'7.647"
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('D41').Value = "'0.5658"
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('D42').Value = "'0.000003098"
$ws.Range('E42').Value = '  +11.33%  '
$ws.Range('D43').Value = "'0.1782"
$ws.Range('E43').Value = '  -1.23%  '
$ws.Range('D44').Value = "'9.400"
$ws.Range('E44').Value = '  -3.13%  '
$ws.Range('D45').Value = "'2.707"
$ws.Range('E45').Value = '  +6.38%  '
$ws.Range('D46').Value = "'0.5266"
$ws.Range('E46').Value = '  -1.47%  '
$ws.Range('D47').Value = "'11.48"
$ws.Range('E47').Value = '  -2.88%  '
$ws.Range('D48').Value = "'0.06860"
$ws.Range('E48').Value = '  -1.65%  '
$ws.Range('D49').Value = "'2.071"
$ws.Range('E49').Value = '  -4.02%  '
$ws.Range('D50').Value = "'1.811"
$ws.Range('E50').Value = '  -1.43%  '
$ws.Range('D51').Value = "'111.27"
$ws.Range('E51').Value = '  -1.53%  '
